$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69; existing rows 69..90 shift down to 70..91.
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new data record.
$ws.Range("A69").Value = 11
$ws.Range("B69").Value = "Vega Monumental Concepción"
$ws.Range("C69").Value = "Bíobío"
$ws.Range("D69").Value = 44463
$ws.Range("E69").Value = 8
$ws.Range("F69").Value = 100112003
$ws.Range("G69").Value = "Ajo"
$ws.Range("H69").Value = "Chino"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 400
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 15500
$ws.Range("M69").Value = 15250
$ws.Range("N69").Value = "$/caja 10 kilos"
$ws.Range("O69").Value = "China"
$ws.Range("P69").Value = 1525
$ws.Range("Q69").Value = 10
$ws.Range("R69").Value = "Hortaliza"
